# Updated cryptos list with GitHub Actions
# Applies the per-cell Price/Volume(1h) refresh (and the Maker/VeChain row swap
# at rows 43-44) produced by the scheduled scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while keeping it text (the sheet stores every Price/Coin/
# Link cell as plain text). Cells whose new text reads as a bare number would
# otherwise get auto-converted to a numeric value on assignment, same as typing
# it into Excel, so those are round-tripped through a Text number format.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    if ($value -match "^[+-]?\d+(\.\d+)?$") {
        $rng.NumberFormat = "@"
        $rng.Value = $value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $value
    }
}

Set-TextValue "D2" "63.476.90"
Set-TextValue "E2" "  +0.48%  "
Set-TextValue "D3" "3.096.69"
Set-TextValue "E3" "  -0.58%  "
Set-TextValue "D5" "583.29"
Set-TextValue "E5" "  -0.25%  "
Set-TextValue "D6" "145.76"
Set-TextValue "E6" "  +0.74%  "
Set-TextValue "E7" "  +0.03%  "
Set-TextValue "D8" "3.090.47"
Set-TextValue "E8" "  -0.58%  "
Set-TextValue "D9" "0.528"
Set-TextValue "E9" "  -0.29%  "
Set-TextValue "D10" "0.161"
Set-TextValue "E10" "  +7.57%  "
Set-TextValue "D11" "5.62"
Set-TextValue "E11" "  -2.34%  "
Set-TextValue "D12" "0.456"
Set-TextValue "E12" "  -2.35%  "
Set-TextValue "E13" "  -1.15%  "
Set-TextValue "D14" "37.24"
Set-TextValue "E14" "  +4.66%  "
Set-TextValue "E15" "  -1.22%  "
Set-TextValue "D16" "3.609.18"
Set-TextValue "E16" "  -0.59%  "
Set-TextValue "D17" "63.324.42"
Set-TextValue "E17" "  +0.40%  "
Set-TextValue "E18" "  -1.16%  "
Set-TextValue "D19" "3.095.22"
Set-TextValue "E19" "  -0.68%  "
Set-TextValue "D20" "461.47"
Set-TextValue "E20" "  -1.31%  "
Set-TextValue "D21" "14.22"
Set-TextValue "E21" "  +1.14%  "
Set-TextValue "D22" "0.722"
Set-TextValue "E22" "  -0.53%  "
Set-TextValue "E23" "  -1.47%  "
Set-TextValue "D24" "81.24"
Set-TextValue "E24" "  -0.89%  "
Set-TextValue "D25" "12.83"
Set-TextValue "E25" "  -3.45%  "
Set-TextValue "E26" "  -1.27%  "
Set-TextValue "E27" "  -0.04%  "
Set-TextValue "D28" "9.00"
Set-TextValue "E28" "  +9.03%  "
Set-TextValue "E29" "  -0.02%  "
Set-TextValue "E30" "  -0.52%  "
Set-TextValue "E31" "  -1.72%  "
Set-TextValue "D32" "6.92"
Set-TextValue "E32" "  +1.21%  "
Set-TextValue "D33" "0.111"
Set-TextValue "E33" "  +0.06%  "
Set-TextValue "D34" "26.62"
Set-TextValue "E34" "  -1.36%  "
Set-TextValue "D35" "0.0₃0846"
Set-TextValue "E35" "  -2.60%  "
Set-TextValue "E36" "  +3.82%  "
Set-TextValue "E37" "  -1.35%  "
Set-TextValue "E38" "  -3.97%  "
Set-TextValue "E39" "  -1.06%  "
Set-TextValue "D40" "50.23"
Set-TextValue "E40" "  -1.31%  "
Set-TextValue "D41" "432.73"
Set-TextValue "E41" "  +0.00%  "
Set-TextValue "D42" "8.67"
Set-TextValue "E42" "  -0.53%  "
Set-TextValue "B43" "Maker"
Set-TextValue "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D43" "2.876.50"
Set-TextValue "E43" "  -1.73%  "
Set-TextValue "B44" "VeChain"
Set-TextValue "C44" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D44" "0.0366"
Set-TextValue "E44" "  -0.89%  "
Set-TextValue "E45" "  -3.51%  "
Set-TextValue "E46" "  -3.39%  "
Set-TextValue "D47" "35.51"
Set-TextValue "E47" "  +1.05%  "
Set-TextValue "D49" "124.06"
Set-TextValue "E49" "  -0.09%  "
Set-TextValue "E50" "  -1.11%  "
Set-TextValue "D51" "24.00"
Set-TextValue "E51" "  -2.03%  "
